# Add a new "type" column (K) to the meta-analysis datasets sheet,
# classifying each dataset as scripted / spontaneous / constrained.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (K1): copy the format of the existing header cell J1
# (bold font, thin left/right border, left/top aligned, wrap text) so the
# new header matches its neighbours, then give it its own distinct style.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "type"

# --- Data column values (K2:K13) ---
$types = @{
    2  = "scripted"
    3  = "scripted"
    4  = "spontaneous"
    5  = "constrained"
    6  = "constrained"
    7  = "constrained"
    8  = "constrained"
    9  = "constrained"
    10 = "scripted"
    11 = "scripted"
    12 = "scripted"
    13 = "constrained"
}

foreach ($row in $types.Keys) {
    $ws.Cells.Item($row, 11).Value = $types[$row]
}

# --- Column width for the new column K ---
$ws.Range("K1").EntireColumn.ColumnWidth = 19

# --- Sheet view tweaks (zoom / scroll position / selection) ---
$window = $excel.ActiveWindow
$window.Zoom = 125
$ws.Range("A9").Select()
$window.ScrollRow = 9
$ws.Range("L11").Select()

Write-Host "done"
